$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the style of A6:A10 back to the plain date style (same as A2:A5),
# using copy/paste-special so the format is deduplicated against the
# existing style index instead of minting a new one.
[void]$ws.Range("A2").Copy()
[void]$ws.Range("A6:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in row 10: new date + new activity text (new shared string)
$ws.Range("A10").Value = 41124
$ws.Range("B10").Value = "Implemented ParallelBitonicSortB8, ParallelBitonicSortB16"

# Update the view: scroll back to A1 and move the selection up to B12:B13
[void]$ws.Range("A1").Select()
[void]$ws.Range("B12:B13").Select()
